$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the query name, index code and condition text in row 2
$ws.Range("A2").Value = "FLPAS_jams"
$ws.Range("B2").Value = "code"
$ws.Range("C2").Value = "code==FLPAS  &  productgroup==10   "

# Narrow column A from ~31.58 to ~24.45 (OOXML char-width units include a
# fixed padding offset of 5/6 of a character versus the COM ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 24.45 - 0.8333333333333334

# Move the view back to the top-left and reselect a cell further down
$ws.Range("A1").Select() | Out-Null
$ws.Range("A8").Select() | Out-Null
